# The commit removes the "KF-N" / "KF-C" columns (M:N) from Sheet1.
# These were two mostly-empty helper columns whose header labels
# ("KF-N"/"KF-C") lived in row 5; deleting the whole columns shifts the
# following "Run" column (formerly O) left into the new column M, and
# Excel automatically drops the now-unused "KF-N"/"KF-C" shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Select columns M:N first so the resulting selection matches what Excel
# leaves behind after an entire-column delete (selection collapses onto
# the cells that slid into the deleted columns' place).
$ws.Columns("M:N").Select()
$ws.Columns("M:N").Delete()
